$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 10: replace the "low" placeholders in B10:D10 with real numbers ---
$ws.Range("B10").Value = 74.2
$ws.Range("C10").Value = 76.8
$ws.Range("D10").Value = 84.4

# --- Row 15: fill in accuracy numbers that were previously blank, bump dropout ---
$ws.Range("B15").Value = 87.3
$ws.Range("C15").Value = 88.9
$ws.Range("D15").Value = 98.9
$ws.Range("O15").Value = 0.7

# --- New row 16 (Trial 15) ---
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 92.4
$ws.Range("C16").Value = 92.3
$ws.Range("D16").Value = 99.9
$ws.Range("E16").Value = "5,5,3,6"
$ws.Range("F16").Value = "2,2"
$ws.Range("G16").Value = "5,5,6,16"
$ws.Range("H16").Value = "2,2"
$ws.Range("K16").Value = 400120
$ws.Range("K16").NumberFormat = "#,##0"
$ws.Range("L16").Value = "120,84"
$ws.Range("M16").Value = "84,43"
$ws.Range("N16").Value = 0.001
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 512
$ws.Range("Q16").Value = 40
$ws.Range("R16").Value = "No"
$ws.Range("S16").Value = "Yes"
$ws.Range("T16").Value = "Generator, No Balancing, 50% augmentation"
$ws.Range("U16").Value = "Xavier Uniform"
$ws.Range("V16").Value = "Zeros"

# --- New row 17 (Trial 16) ---
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = 94
$ws.Range("C17").Value = 94.5
$ws.Range("D17").Value = 99.8
$ws.Range("E17").Value = "5,5,3,9"
$ws.Range("F17").Value = "2,2"
$ws.Range("G17").Value = "5,5,6,20"
$ws.Range("H17").Value = "2,2"
$ws.Range("K17").Value = 500500
$ws.Range("L17").Value = 500250
$ws.Range("K17:L17").NumberFormat = "#,##0"
$ws.Range("M17").Value = "250,43"
$ws.Range("N17").Value = 0.001
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 512
$ws.Range("Q17").Value = 40
$ws.Range("R17").Value = "No"
$ws.Range("S17").Value = "Yes"
$ws.Range("T17").Value = "Generator, No Balancing, 50% augmentation"
$ws.Range("U17").Value = "Xavier Uniform"
$ws.Range("V17").Value = "Zeros"

# --- New row 18 (Trial 17 - the best run, highlighted green) ---
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 94.6
$ws.Range("C18").Value = 95.6
$ws.Range("D18").Value = 99.6
$ws.Range("E18").Value = "5,5,3,9"
$ws.Range("F18").Value = "2,2"
$ws.Range("G18").Value = "5,5,6,18"
$ws.Range("I18").Value = "5,5,18,27"
$ws.Range("J18").Value = "2,2"
$ws.Range("K18").Value = 243243
$ws.Range("L18").Value = 243243
$ws.Range("M18").Value = "243,43"
$ws.Range("N18").Value = 0.001
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 512
$ws.Range("Q18").Value = 40
$ws.Range("R18").Value = "No"
$ws.Range("S18").Value = "Yes"
$ws.Range("T18").Value = "Generator, No Balancing, 50% augmentation"
$ws.Range("U18").Value = "Xavier Uniform"
$ws.Range("V18").Value = "Zeros"

# Highlight the final/best trial row with a green fill + explicit Calibri font,
# and keep the thousands-separator format on the FC1/FC2 columns.
$row18 = $ws.Range("A18:V18")
$row18.Font.Name = "Calibri"
$row18.Interior.Color = 5287936
$ws.Range("K18:L18").NumberFormat = "#,##0"

# --- Sheet-level cosmetics ---
$ws.PageSetup.Orientation = 1
$ws.Range("D10").Select()
